# Generate Report for Handback
# Renames the prior "f11706e4-..." handback record to "9617c70a-..." (in
# place, simulating a re-run against the same slot) and appends a brand new
# record for "c2ca4700-...".

$wb = $excel.ActiveWorkbook

$oldGuid = "f11706e4-eeb2-43d4-bd78-9291a933d500"
$newGuid = "9617c70a-46f5-46ac-ab30-2f9d41ca007b"
$addGuid = "c2ca4700-64fb-4643-8339-4a0566cf22e3"

$genDate        = "2016-08-15 12:56:13"
$zhHandoffDate  = "2016-08-15 12:56:07"
$zhHandbackDate = "2016-08-15 12:56:33"
$deHandbackDate = "2016-08-15 12:56:39"

$zhHashNew = "a273a745e1ad84d00ce66e295e2e2c2dc027e5b1"
$zhHashAdd = "17ee96effe9bd37e690314f60cec256807654bbe"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)

# -- update the existing row (row 2): rename the file + bump the date
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0c271e8fc35a272663451828cb9359618eb1833a/e2e/$newGuid.md",
    "",
    "",
    "e2e\$newGuid.md") | Out-Null
$wsOverview.Range("G2").Value = $genDate

# -- append new row (row 3) for the newly handed-back file
$loOverview.ListRows.Add() | Out-Null
$wsOverview.Range("A3").Value = "$addGuid.md"
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0c271e8fc35a272663451828cb9359618eb1833a/e2e/$addGuid.md",
    "",
    "",
    "e2e\$addGuid.md") | Out-Null
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G3").Value = $genDate

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)

# -- update the existing row (row 2)
$wsZh.Range("A2").Hyperlinks.Delete()
$wsZh.Hyperlinks.Add(
    $wsZh.Range("A2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0c271e8fc35a272663451828cb9359618eb1833a/e2e/$newGuid.md",
    "",
    "",
    "$newGuid.md") | Out-Null
$wsZh.Range("G2").Value = "$newGuid.$zhHashNew.zh-cn.xlf"
$wsZh.Range("H2").Value = $zhHandoffDate
$wsZh.Range("I2").Hyperlinks.Delete()
$wsZh.Hyperlinks.Add(
    $wsZh.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/beeaa8091dd05d4a2e89a3c392df0c24033b0957/e2e/$newGuid.md",
    "",
    "",
    "$newGuid.md") | Out-Null
$wsZh.Range("J2").Value = "$newGuid.$zhHashNew.zh-cn.xlf"
$wsZh.Range("K2").Value = $zhHandbackDate

# -- append new row (row 3)
$loZh.ListRows.Add() | Out-Null
$wsZh.Hyperlinks.Add(
    $wsZh.Range("A3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0c271e8fc35a272663451828cb9359618eb1833a/e2e/$addGuid.md",
    "",
    "",
    "$addGuid.md") | Out-Null
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "'True"
$wsZh.Range("G3").Value = "$addGuid.$zhHashAdd.zh-cn.xlf"
$wsZh.Range("H3").Value = $zhHandoffDate
$wsZh.Hyperlinks.Add(
    $wsZh.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/beeaa8091dd05d4a2e89a3c392df0c24033b0957/e2e/$addGuid.md",
    "",
    "",
    "$addGuid.md") | Out-Null
$wsZh.Range("J3").Value = "$addGuid.$zhHashAdd.zh-cn.xlf"
$wsZh.Range("K3").Value = $zhHandbackDate
$wsZh.Range("L3").Value = ""
$wsZh.Range("M3").Value = "'True"
$wsZh.Range("N3").Value = ""
$wsZh.Range("O3").Value = "'False"
$wsZh.Range("P3").Value = ""

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)

# -- update the existing row (row 2)
$wsDe.Range("A2").Hyperlinks.Delete()
$wsDe.Hyperlinks.Add(
    $wsDe.Range("A2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0c271e8fc35a272663451828cb9359618eb1833a/e2e/$newGuid.md",
    "",
    "",
    "$newGuid.md") | Out-Null
$wsDe.Range("G2").Value = "$newGuid.$zhHashNew.de-de.xlf"
$wsDe.Range("H2").Value = $genDate
$wsDe.Range("I2").Hyperlinks.Delete()
$wsDe.Hyperlinks.Add(
    $wsDe.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/c76ca424be44f43384e106c34c0ff6eadc196c84/e2e/$newGuid.md",
    "",
    "",
    "$newGuid.md") | Out-Null
$wsDe.Range("J2").Value = "$newGuid.$zhHashNew.de-de.xlf"
$wsDe.Range("K2").Value = $deHandbackDate

# -- append new row (row 3)
$loDe.ListRows.Add() | Out-Null
$wsDe.Hyperlinks.Add(
    $wsDe.Range("A3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0c271e8fc35a272663451828cb9359618eb1833a/e2e/$addGuid.md",
    "",
    "",
    "$addGuid.md") | Out-Null
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "'True"
$wsDe.Range("G3").Value = "$addGuid.$zhHashAdd.de-de.xlf"
$wsDe.Range("H3").Value = $genDate
$wsDe.Hyperlinks.Add(
    $wsDe.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/c76ca424be44f43384e106c34c0ff6eadc196c84/e2e/$addGuid.md",
    "",
    "",
    "$addGuid.md") | Out-Null
$wsDe.Range("J3").Value = "$addGuid.$zhHashAdd.de-de.xlf"
$wsDe.Range("K3").Value = $deHandbackDate
$wsDe.Range("L3").Value = ""
$wsDe.Range("M3").Value = "'True"
$wsDe.Range("N3").Value = ""
$wsDe.Range("O3").Value = "'False"
$wsDe.Range("P3").Value = ""
